# Add separability results in csv
# Round the ConvexHullArea values (column D) to the nearest whole number
# for all data rows (rows 2-198) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 198
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $v = [double]$cell.Value2
    $cell.Value = [Math]::Round($v, 0)
}
